$wb = $excel.ActiveWorkbook
$originalSheet = $wb.ActiveSheet
$ws = $wb.Worksheets.Item("EDLE")

# Update the dispatch logit exponent value for "all electricity sources"
$ws.Range("B2").Value = -1

# Update the saved selection on the EDLE sheet to B3
$ws.Activate()
$ws.Range("B3").Select()

# Restore whichever sheet/tab was active before this edit
$originalSheet.Activate()
